# Add AX5043 and BFP740F LNA BOMs
#
# - Boards sheet: mark the BFP740F LNA (row 6) and AX5043 transceiver (row 24)
#   boards as received, record quantities, and bump their BOM status from
#   PEND to ORDRD. Also flips the TRL cal board (row 25) and the BGB741
#   LNA BOM (row 7) to received, and renames the quad hybrid part.
# - Inventory sheet: two blank history rows added near the bottom, and a
#   resistor count ticked down by one (16 -> 15).

$wb = $excel.ActiveWorkbook

$boards = $wb.Worksheets.Item("Boards")
$inventory = $wb.Worksheets.Item("Inventory")

# Color constants (Excel COM uses 0xBBGGRR ordering for RGB()-style values)
$colorGreen  = 5296274   # fgColor FF92D050 - "RCVD" rows (style index 15)
$colorOrange = 49407     # fgColor FFFFC000 - "ORDRD" rows (style index 20)

# --- Boards!row6 : amp-LNA-BFP740Fxxx -------------------------------------
$boards.Range("B6").Value = "RCVD"
$boards.Range("B6").Interior.Color = $colorGreen
$boards.Range("C6").Value = 3
$boards.Range("E6").Value = "ORDRD"
$boards.Range("E6").Interior.Color = $colorOrange

# --- Boards!row7 : amp-LNA-BGB741L7ESD ------------------------------------
$boards.Range("E7").Value = "RCVD"
$boards.Range("E7").Interior.Color = $colorGreen

# --- Boards!row21 : quad-hybrid-QCN-19D+ -> quad-hybrid-QCN-xxx+ ---------
$boards.Range("A21").Value = "quad-hybrid-QCN-xxx+"
$boards.Range("B21").Value = "ORDRD"
$boards.Range("B21").Interior.Color = $colorOrange
$boards.Range("E21").Value = "ORDRD"
$boards.Range("E21").Interior.Color = $colorOrange

# --- Boards!row24 : transceiver-AX5043 ------------------------------------
$boards.Range("B24").Value = "RCVD"
$boards.Range("B24").Interior.Color = $colorGreen
$boards.Range("C24").Value = 3
$boards.Range("E24").Value = "ORDRD"
$boards.Range("E24").Interior.Color = $colorOrange

# --- Boards!row25 : TRL-calibration-thru-reflect-match --------------------
$boards.Range("B25").Value = "RCVD"
$boards.Range("B25").Interior.Color = $colorGreen
$boards.Range("C25").Value = 1
$boards.Range("D25").Value = 2

# --- Inventory sheet -------------------------------------------------------
# Two new blank history rows, inserted right before the final (bordered)
# history row so the bottom border stays on the last row.
$inventory.Rows("100:101").Insert()

# Resistor count decremented by one.
$inventory.Range("A53").Value = 15
